$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "26.330.96" }
    @{ Cell = "E2"; Value = "  -3.06%  " }
    @{ Cell = "D3"; Value = "1.790.84" }
    @{ Cell = "E3"; Value = "  -2.97%  " }
    @{ Cell = "E4"; Value = "  +0.53%  " }
    @{ Cell = "D5"; Value = "1.007" }
    @{ Cell = "E5"; Value = "  +0.49%  " }
    @{ Cell = "D6"; Value = "306.51" }
    @{ Cell = "E6"; Value = "  -2.18%  " }
    @{ Cell = "D7"; Value = "0.4547" }
    @{ Cell = "E7"; Value = "  -1.72%  " }
    @{ Cell = "D8"; Value = "0.3620" }
    @{ Cell = "E8"; Value = "  -1.90%  " }
    @{ Cell = "E9"; Value = "  -2.73%  " }
    @{ Cell = "D10"; Value = "0.8704" }
    @{ Cell = "E10"; Value = "  -1.60%  " }
    @{ Cell = "D11"; Value = "0.07785" }
    @{ Cell = "E11"; Value = "  -0.38%  " }
    @{ Cell = "D12"; Value = "19.37" }
    @{ Cell = "E12"; Value = "  -2.47%  " }
    @{ Cell = "D13"; Value = "1.789.88" }
    @{ Cell = "E13"; Value = "  -3.94%  " }
    @{ Cell = "D14"; Value = "5.254" }
    @{ Cell = "E14"; Value = "  -2.48%  " }
    @{ Cell = "D15"; Value = "6.304" }
    @{ Cell = "E15"; Value = "  -2.91%  " }
    @{ Cell = "D16"; Value = "84.33" }
    @{ Cell = "E16"; Value = "  -7.63%  " }
    @{ Cell = "E17"; Value = "  +0.66%  " }
    @{ Cell = "D18"; Value = "0.000008486" }
    @{ Cell = "E18"; Value = "  -3.94%  " }
    @{ Cell = "E19"; Value = "  +0.44%  " }
    @{ Cell = "D20"; Value = "26.396.47" }
    @{ Cell = "E20"; Value = "  -2.92%  " }
    @{ Cell = "D21"; Value = "14.14" }
    @{ Cell = "E21"; Value = "  -3.28%  " }
    @{ Cell = "D22"; Value = "4.966" }
    @{ Cell = "E22"; Value = "  -1.38%  " }
    @{ Cell = "E23"; Value = "  -0.46%  " }
    @{ Cell = "D24"; Value = "2.001.72" }
    @{ Cell = "E24"; Value = "  -5.27%  " }
    @{ Cell = "D25"; Value = "1.979" }
    @{ Cell = "E25"; Value = "  -2.62%  " }
    @{ Cell = "D26"; Value = "152.07" }
    @{ Cell = "E26"; Value = "  +0.72%  " }
    @{ Cell = "D27"; Value = "17.79" }
    @{ Cell = "E27"; Value = "  -3.02%  " }
    @{ Cell = "D28"; Value = "2.035" }
    @{ Cell = "E28"; Value = "  +0.70%  " }
    @{ Cell = "D29"; Value = "112.17" }
    @{ Cell = "E29"; Value = "  -2.93%  " }
    @{ Cell = "D30"; Value = "4.821" }
    @{ Cell = "E30"; Value = "  -3.61%  " }
    @{ Cell = "D31"; Value = "0.08654" }
    @{ Cell = "E31"; Value = "  -2.26%  " }
    @{ Cell = "E32"; Value = "  -4.19%  " }
    @{ Cell = "D33"; Value = "4.430" }
    @{ Cell = "E33"; Value = "  -1.55%  " }
    @{ Cell = "B34"; Value = "ImmutableX" }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D34"; Value = "0.7132" }
    @{ Cell = "E34"; Value = "  -8.67%  " }
    @{ Cell = "B35"; Value = "RenderToken" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = "D35"; Value = "2.653" }
    @{ Cell = "E35"; Value = "  -1.76%  " }
    @{ Cell = "E36"; Value = "  -3.58%  " }
    @{ Cell = "D37"; Value = "1.006" }
    @{ Cell = "E37"; Value = "  +0.25%  " }
    @{ Cell = "D38"; Value = "1.078" }
    @{ Cell = "E38"; Value = "  -2.25%  " }
    @{ Cell = "D39"; Value = "0.01935" }
    @{ Cell = "E39"; Value = "  -0.37%  " }
    @{ Cell = "D40"; Value = "0.05079" }
    @{ Cell = "E40"; Value = "  -2.43%  " }
    @{ Cell = "D41"; Value = "2.859" }
    @{ Cell = "E41"; Value = "  -3.15%  " }
    @{ Cell = "D42"; Value = "6.884" }
    @{ Cell = "E42"; Value = "  -1.81%  " }
    @{ Cell = "D43"; Value = "0.4907" }
    @{ Cell = "E43"; Value = "  -2.40%  " }
    @{ Cell = "D44"; Value = "0.1511" }
    @{ Cell = "E44"; Value = "  -6.10%  " }
    @{ Cell = "D45"; Value = "7.959" }
    @{ Cell = "E45"; Value = "  -5.94%  " }
    @{ Cell = "D46"; Value = "1.008" }
    @{ Cell = "D47"; Value = "0.4557" }
    @{ Cell = "E47"; Value = "  -3.88%  " }
    @{ Cell = "D48"; Value = "9.846" }
    @{ Cell = "E48"; Value = "  -4.76%  " }
    @{ Cell = "D49"; Value = "99.58" }
    @{ Cell = "E49"; Value = "  -3.20%  " }
    @{ Cell = "D50"; Value = "1.578" }
    @{ Cell = "E50"; Value = "  -3.14%  " }
    @{ Cell = "D51"; Value = "0.05964" }
    @{ Cell = "E51"; Value = "  -3.61%  " }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
